$d = $word.ActiveDocument

# Rename the bookmark from "bookmark_test" to "ABCD-1234" by re-adding it
# over the same range (Word has no direct bookmark-rename API).
$bm = $d.Bookmarks.Item("bookmark_test")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("ABCD-1234", $bmRange)

# Change the bookmarked word "Bookmark" -> "Bookmarked".
$null = $d.Content.Find.Execute("Bookmark", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Bookmarked", 2)

# Drop the trailing " Hello" text that followed the bookmark.
$null = $d.Content.Find.Execute(" Hello", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 2)

# Remove the whole second paragraph ("World"), including its paragraph mark.
$d.Paragraphs.Item(2).Range.Delete()
